$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(' Kumele uwacedze onkhe emakhono langaphansi kwemgomo kute utfole i-mbheji. Nase uwucedzile munye umgomo, ungakhetsa lomunye. Uma sewutfole onkhe emabheji akho, utawube sewucedzile kufundza futsi utawutfola umklomelo wekukhulisa bantfwana ngendlela lekahle.', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ' Kumele uwacedze onkhe emakhono langaphansi kwemgomo kute utfole i-mbheji. Nase uwucedzile munye umgomo, ungakhetsa lomunye. Nase utitfolile tonkhe timbheji talemigomo leyehlukene, loko kuso kutsi utawube sewuticedzile tifundvo takho futsi utawutfola umklomelo wekukhulisa bantfwana ngendlela lekahle.' } else { Write-Host 'NOT FOUND [0]:  Kumele uwacedze onkhe emakhono langaphansi kwemgo' }
$rng = $d.Content
$found = $rng.Find.Execute('Njengobe uchubeka nesifundvo ngasinye, utawutfola lwati lolubonakala ngalendlela. Letibonakaliso tikutjela kutsi sewuhambile kangakanani esifundvweni. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Nawusachubeka netifundvo onkhe malanga, utawutfola umlayeto lokwatisa ngendlela lochuba ngayo, lotawubukeka kanjena; . Letimphawu tikutjela kutsi sewuhambe kangakanani esifundvweni sakho sangalelo langa. ' } else { Write-Host 'NOT FOUND [1]: Njengobe uchubeka nesifundvo ngasinye, utawutfola ' }
$rng = $d.Content
$found = $rng.Find.Execute(' Nangabe ufuna kubona kutsi sewutfutfuke kanganani ekufinyeleleni imigomo yakho yekuba ngumtali, ungahlola inchubekelembili loyitfolile ngekusebentisa i-Main Menyu. Kute ufinyelele imenyu, bhala "Imenyu" nobe nini. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ' Nawufuna kubona kutsi sewutfutfuke kanganani ekufinyeleleni kulemigomo yakho yekuba ngumtali lokahle, ungahlola indlela lewuchuba ngayo ku-Main Menyu. Kute ufinyelele ku-menyu, bhala "Imenyu" noma ngabe ngunini ' } else { Write-Host 'NOT FOUND [2]:  Nangabe ufuna kubona kutsi sewutfutfuke kanganani' }
$rng = $d.Content
$found = $rng.Find.Execute(' Inketfo yekucala ngu-menyu ibhalwe "Landzelela inchubekelembili yami". ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ' Lotakukhetsa kucala kubhalwe kutsi "landzelela indlela lengichuba ngayo". ' } else { Write-Host 'NOT FOUND [3]:  Inketfo yekucala ngu-menyu ibhalwe "Landzelela in' }
$rng = $d.Content
$found = $rng.Find.Execute(' Lapha ungabona inchubekelembili yakho, uphindze ubuke imigomo loyifinyelele kanye naleyo lengakacedvwa.', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ' Lapha ungabona indlela lochuba ngayo, ubukete imigomo lose uyicedzile naleyo lekusemele kutsi uyente.' } else { Write-Host 'NOT FOUND [4]:  Lapha ungabona inchubekelembili yakho, uphindze u' }
$rng = $d.Content
$found = $rng.Find.Execute('Lemenyu inetintfo letingakusita. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Le-menyu icuketse letinye tintfo letingakusita. ' } else { Write-Host 'NOT FOUND [5]: Lemenyu inetintfo letingakusita. ' }
$rng = $d.Content
$found = $rng.Find.Execute('Ngetulu kwekulandzelela inchubekelembili yakho, ungaphindze: ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Ngetulu kwekulandzelela indlela lochuba ngayo, ungaphindze: ' } else { Write-Host 'NOT FOUND [6]: Ngetulu kwekulandzelela inchubekelembili yakho, un' }
$rng = $d.Content
$found = $rng.Find.Execute('Cocelenani nge-ParentText nemngani wakho futsi umsite abhalise. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Cocela umngani wakho nge-ParentText bese uyamsita kutsi abhalise. ' } else { Write-Host 'NOT FOUND [7]: Cocelenani nge-ParentText nemngani wakho futsi ums' }
$rng = $d.Content
$found = $rng.Find.Execute('Shintja amasethingi akho, njengendlela lowatfola ngayo umlayeto, sikhatsi lowatiswa ngaso, nobe ulungise lwati ngawe kanye newemntfwana wakho kute utfole lusito lolufanele. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Shintja amasethingi akho, njengendlela lowatfola ngayo umlayeto, sikhatsi sekutfola satiso, noma ulungise imininingwane yakho kanye nemntfwana wakho kute nitfole lusito lolufanele. ' } else { Write-Host 'NOT FOUND [8]: Shintja amasethingi akho, njengendlela lowatfola n' }
$rng = $d.Content
$found = $rng.Find.Execute('Bhala luhla lwetintfo longatenta nemntfwanakho kute wakhe buhlobo benu. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Tfola luhla lwetintfo longatenta nemntfwana wakho kute nakhe buhlobo benu. ' } else { Write-Host 'NOT FOUND [9]: Bhala luhla lwetintfo longatenta nemntfwanakho kut' }
$rng = $d.Content
$found = $rng.Find.Execute('Hlola sicondziso sekungena bese utfola lusito lwengusebentisa iParentText. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Buketa lesicondziso sekungenela loluhlelo bese utfola lusito lwekusebentisa iParentText. ' } else { Write-Host 'NOT FOUND [10]: Hlola sicondziso sekungena bese utfola lusito lwen' }
$rng = $d.Content
$found = $rng.Find.Execute('Futsi tfole lusito lwekusombulula tinkinga letimatima umntfwana wakho lahlangabetana nato. Asesifundze kabanti ngalendzaba nyalo. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Futsi utfole lusito lwekusombulula tinkinga letimatima lakahlangabetana nato umntfwana wakho. Asesifundze kabanti ngalesici nyalo. ' } else { Write-Host 'NOT FOUND [11]: Futsi tfole lusito lwekusombulula tinkinga letimat' }
$rng = $d.Content
$found = $rng.Find.Execute('Kuba ngumtali kungaba matima. Kuba ngumtali kungaba matima. Ngisho nobe tinkinga utiva tihlukile kuwe, kodvwa tivame kakhulu kunaloko locabanga kutsi tiyenteka. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Kuba ngumtali kungaba matima. Kungenteka ucabange kutsi tinkinga tikhungetse wena wedva, tivame kakhulu kunaloku wena lokucabangako. ' } else { Write-Host 'NOT FOUND [12]: Kuba ngumtali kungaba matima. Kuba ngumtali kungab' }
$rng = $d.Content
$found = $rng.Find.Execute('Njengobe ucala kufinyelela imigomo yakho kuloluhlelo, ngitawubuya ngitewubuta kutsi kuhamba njani umntfwana wakho. Ngingase ngibanikete lusito. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Nase ucala kuyifeza lemigomo yakho kuletifundvo, ngitawubuya ngitohlola kutsi nichuba kanjani kanye nemntfwana wakho. Nangabe kukhona lokungahambi kahle, ngingahle nginisite. ' } else { Write-Host 'NOT FOUND [13]: Njengobe ucala kufinyelela imigomo yakho kuloluhle' }
$rng = $d.Content
$found = $rng.Find.Execute('Nangabe nicoca nami ngetinkinga leninato, ngitaniniketa tisombululo letitawusita niphumelele. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Uma ningicocela ngetinkinga leninato, ngitaniniketa tisombululo letitawunisita niphumelele. ' } else { Write-Host 'NOT FOUND [14]: Nangabe nicoca nami ngetinkinga leninato, ngitanin' }
$rng = $d.Content
$found = $rng.Find.Execute('Kodvwa, akudzingeki ulindzele mina kutsi ngikusekele. Ungaphindze utfole lusito lwekulungisa tinkinga ngu-Main Menu nobe nini. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Kodvwa loko akusho kutsi bese nimela mine kutsi ngininike lusito,. Ungaphindze utfole lusito lwekulungisa tinkinga ku-Main Menyu noma nini. ' } else { Write-Host 'NOT FOUND [15]: Kodvwa, akudzingeki ulindzele mina kutsi ngikuseke' }
$rng = $d.Content
$found = $rng.Find.Execute('Lusito lwetimo letiphutfumako kanye Netinkinga', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Lusito lwetimo letiphutfumako kanye Netinhlekelele' } else { Write-Host 'NOT FOUND [16]: Lusito lwetimo letiphutfumako kanye Netinkinga' }
$rng = $d.Content
$found = $rng.Find.Execute('Nangabe udzinga lwati ngetinsita emmangweni wakini tekubhekana nebudlova basemndenini, budlova lobentiwa ngekwelicasi, kuphatfwa kabi kwengcondvo, nobe letinye timo letiphutfumako, ungabhala umlayeto ku-LUSITO nobe kunini futsi utfole imininingwane yekutsintsana nebantfu labangakusita. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Nangabe udzinga lwati ngetinsita emmangweni wakini tekubhekana nebudlova basemndenini, budlova lobentiwa ngekwelicasi, kuphatseka kabi kwengcondvo, noma letinye timo letiphutfumako, ungabhala umlayeto lotsi LUSITO noma ngabe kunini kute utfole imininingwane yebantfu labangakusita. ' } else { Write-Host 'NOT FOUND [17]: Nangabe udzinga lwati ngetinsita emmangweni wakini' }
$rng = $d.Content
$found = $rng.Find.Execute('Lwati lwakho luvikelekile: Kute lutfo lolutawudluliselwa ngaphandle kwemvume yakho futsi ngeke kutsengiswe ngalo kute kuzuze. Umlayeto lowutfumelako ubhaliwe futsi ugcinwe endzaweni levikelekile. Nobe kunjalo, khumbula kutsi nangabe lotsite atfola lucingo lwakho futsi aluvule, angakhona kubuka imilayeto yakho kute abone loko bhalile. Nangabe utfumela umniningwane lobalulekile futsi ukhatsatekile, ciniseka kutsi uyawususa lomlayeto elucingweni lwakho. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Imininingwane yakho iphephile lapha: angeke yendluliselwe kulabanye ngaphandle kwemvumo yakho futsi angeke kutsengiswe ngayo kute kutfolakale inzuzo. Umlayeto lowutfumelako ubhalwe ngendlela lebangeke bakhone kuyifundza labanye futsi ugcinwe endzaweni levikelekile. Noma kunjalo, khumbula kutsi nangabe umuntfu lotsite atfola lucingo lwakho futsi aluvule, angakhona kubuka imilayeto yakho kute abone loko lokubhalile. Nangabe utfumela umlayeto lobucayi ukhatsatekile, bani nesiciniseko kutsi uyawucisha lomlayeto elucingweni lwakho. ' } else { Write-Host 'NOT FOUND [18]: Lwati lwakho luvikelekile: Kute lutfo lolutawudlul' }
$rng = $d.Content
$found = $rng.Find.Execute('Ngiyabonga kakhulu ngekungilalela! Siyetsemba kutsi utawujabulela luhambo lwakho lweParentText futsi ulisebentise kahle! Ungayi tfola le-video nobe nini ngu-main menyu. ', $true, $true, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = 'Ngiyabonga kakhulu kungilalela! Siyetsemba kutsi utawulujabulela luhambo lwakho lweParentText futsi ulisebentise kahle! Ungayitfola le-video noma kunini ku-main menyu. ' } else { Write-Host 'NOT FOUND [19]: Ngiyabonga kakhulu ngekungilalela! Siyetsemba kuts' }
